$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

$ws.Range("A8").Value = "intro_1"
$ws.Range("B8").Value = "intro 1"
$ws.Range("A9").Value = "intro_2"
$ws.Range("B9").Value = "intro 2"

$ws.Range("A10").Select()
